# 2013-11-22  审批  修复    无法添加
# 2013-11-22  JS    删除    删除google cdn jquery
#
# Append two new log rows (9 and 10) below the existing change-log table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: 审批 / 修复 entry, with related-table info (flow / - / emp_no)
$ws.Range("A9").Value = "2013-11-22"
$ws.Range("B9").Value = "审批"
$ws.Range("C9").Value = "修复"
$ws.Range("E9").Value = "flow"
$ws.Range("F9").Value = "-"
$ws.Range("G9").Value = "emp_no"

# Row 10: JS / 删除 entry - removed google cdn jquery
$ws.Range("A10").Value = "2013-11-22"
$ws.Range("B10").Value = "JS"
$ws.Range("C10").Value = "删除"
$ws.Range("D10").Value = "删除google cdn jquery"

# Match the row height used by every other row in the table
$ws.Rows.Item(9).RowHeight = 21.6
$ws.Rows.Item(10).RowHeight = 21.6

# Leave the selection on the last-edited cell, like a freshly-typed row
$ws.Range("A10").Select() | Out-Null
